# Correcting market share tab for updated scenario 3s
#
# The "New Product A" / "Old Product B (SOC)" hand-off for DRC/MDA was
# pushed out: New Product A now only starts picking up share in 2026
# (column L, since D=2018 ... L=2026 ... Z=2040) instead of 2018, and
# Old Product B (SOC) now loses its share in 2026 instead of holding it
# all the way through 2040.

$wb = $excel.ActiveWorkbook
$platformCoverage = $wb.Worksheets.Item("Platform Coverage")
$marketShare = $wb.Worksheets.Item("MarketShare")

# New Product A (row 2): add 100% share for 2026-2040 (cols L:Z).
$marketShare.Range("L2:Z2").Value = 1

# Old Product B (SOC) (row 3): it now only keeps its 100% share through
# 2025 (cols D:K); the 2026-2040 figures (cols L:Z) move up to row 2.
$marketShare.Range("L3:Z3").ClearContents()

# Reflect the edit in the UI state: the MarketShare tab becomes the
# active/selected sheet, with the newly-entered run of cells selected.
$marketShare.Activate()
$marketShare.Range("L2:Z2").Select()
